# Hide Problem Set 22 slides.
$p = $ppt.ActivePresentation

# Slides 32 through 40 (inclusive) correspond to slides/slide32.xml .. slide40.xml,
# which cover the "PS 22.x" Problem Set 22 recap/tutorial slides. Mark them hidden
# in the slide show (adds show="0" to the <p:sld> element).
for ($i = 32; $i -le 40; $i++) {
    $slide = $p.Slides.Item($i)
    $slide.SlideShowTransition.Hidden = $true
}

# Slide 32 also drops the leading "Recap. " run from the subtitle text box, leaving
# "PS 22.1. PS 22.2" (each still in its own differently-coloured run).
$slide32 = $p.Slides.Item(32)
$subtitle = $slide32.Shapes.Item(4)
$textRange = $subtitle.TextFrame.TextRange
$firstRun = $textRange.Runs(1)
if ($firstRun.Text -eq "Recap. ") {
    $firstRun.Text = ""
}
